$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is the daily price row being updated with the latest automatic data refresh.
$ws.Range("A2").Value = 45932
$ws.Range("B2").Value = 104.34
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 95
$ws.Range("E2").Value = 97.06999999999999
$ws.Range("F2").Value = 100.32
$ws.Range("G2").Value = 105.02
$ws.Range("H2").Value = 105.49
$ws.Range("I2").Value = 121.42
$ws.Range("J2").Value = 131.15
$ws.Range("K2").Value = 93.62
$ws.Range("L2").Value = 59.12
$ws.Range("M2").Value = 34.29
$ws.Range("N2").Value = 22.55
$ws.Range("O2").Value = 17.67
$ws.Range("P2").Value = 19.26
$ws.Range("Q2").Value = 16.59
$ws.Range("R2").Value = 31.9
$ws.Range("S2").Value = 45.64
$ws.Range("T2").Value = 80.18000000000001
$ws.Range("U2").Value = 106.92
$ws.Range("V2").Value = 204.26
$ws.Range("W2").Value = 156.49
$ws.Range("X2").Value = 118.79
$ws.Range("Y2").Value = 104.7
$ws.Range("Z2").Value = 86.31999999999999
$ws.Range("AB2").Value = 146.06
$ws.Range("AD2").Value = 180.38
$ws.Range("AE2").Value = "6h-8h"
$ws.Range("AF2").Value = 113.46

$wb.Save()
